$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-detected as a number by
# the Value setter (e.g. "580.64") are written via a text-number-format round
# trip so they stay text cells, matching the original inlineStr cell type.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '66.535.39'
$ws.Range('E2').Value = '  -1.13%  '

$ws.Range('D3').Value = '3.451.64'
$ws.Range('E3').Value = '  -0.83%  '

$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue 'D5' '580.64'
$ws.Range('E5').Value = '  -2.18%  '

Set-TextValue 'D6' '176.03'
$ws.Range('E6').Value = '  -1.33%  '

$ws.Range('E7').Value = '  -0.01%  '

Set-TextValue 'D8' '0.598'
$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('D9').Value = '3.449.27'
$ws.Range('E9').Value = '  -0.94%  '

$ws.Range('E10').Value = '  -2.36%  '

$ws.Range('E11').Value = '  -3.12%  '

Set-TextValue 'D12' '0.419'
$ws.Range('E12').Value = '  -3.47%  '

$ws.Range('D13').Value = '4.043.88'
$ws.Range('E13').Value = '  -0.89%  '

Set-TextValue 'D14' '30.59'
$ws.Range('E14').Value = '  -4.28%  '

$ws.Range('E15').Value = '  -3.34%  '

$ws.Range('D16').Value = '66.517.29'
$ws.Range('E16').Value = '  -1.23%  '

$ws.Range('E17').Value = '  -2.68%  '

$ws.Range('D18').Value = '3.448.96'
$ws.Range('E18').Value = '  -0.88%  '

Set-TextValue 'D19' '6.01'
$ws.Range('E19').Value = '  -3.92%  '

$ws.Range('E20').Value = '  -3.02%  '

Set-TextValue 'D21' '376.26'
$ws.Range('E21').Value = '  -3.18%  '

$ws.Range('E22').Value = '  -2.50%  '

$ws.Range('E23').Value = '  +0.04%  '

Set-TextValue 'D24' '5.72'
$ws.Range('E24').Value = '  -0.06%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D25' '70.98'
$ws.Range('E25').Value = '  -4.09%  '

$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D26' '0.527'
$ws.Range('E26').Value = '  -1.50%  '

$ws.Range('E27').Value = '  -3.06%  '

Set-TextValue 'D28' '9.82'
$ws.Range('E28').Value = '  -5.13%  '

$ws.Range('E29').Value = '  -1.78%  '

$ws.Range('E30').Value = '  +0.10%  '

Set-TextValue 'D31' '5.86'
$ws.Range('E31').Value = '  -4.92%  '

Set-TextValue 'D32' '23.97'
$ws.Range('E32').Value = '  +1.88%  '

$ws.Range('E33').Value = '  -3.52%  '

$ws.Range('E34').Value = '  -5.28%  '

Set-TextValue 'D35' '0.999'
$ws.Range('E35').Value = '  -0.07%  '

Set-TextValue 'D36' '7.03'
$ws.Range('E36').Value = '  -4.30%  '

$ws.Range('E37').Value = '  -4.80%  '

Set-TextValue 'D38' '159.44'
$ws.Range('E38').Value = '  -2.84%  '

Set-TextValue 'D39' '0.877'
$ws.Range('E39').Value = '  +0.56%  '

Set-TextValue 'D40' '27.28'
$ws.Range('E40').Value = '  +3.95%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D41' '1.78'
$ws.Range('E41').Value = '  -4.90%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D42' '2.63'
$ws.Range('E42').Value = '  -3.40%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D43' '6.48'
$ws.Range('E43').Value = '  -5.56%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D44' '4.46'
$ws.Range('E44').Value = '  -3.66%  '

$ws.Range('D45').Value = '2.685.01'
$ws.Range('E45').Value = '  -5.57%  '

$ws.Range('E46').Value = '  -3.39%  '

Set-TextValue 'D47' '25.23'
$ws.Range('E47').Value = '  -6.08%  '

Set-TextValue 'D48' '40.20'
$ws.Range('E48').Value = '  -3.13%  '

$ws.Range('E49').Value = '  -1.69%  '

Set-TextValue 'D50' '321.09'
$ws.Range('E50').Value = '  -4.42%  '

$ws.Range('E51').Value = '  -4.01%  '

